# Bug fix in init: fixed bug when selecting container type that's allowed.
# Replaces the (mis-parsed) CBM readings with the corrected values and
# extends the list to the full 30-row result set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  21.31,
  89.59,
  9.35,
  68.040000000000006,
  82.91,
  1767.45,
  5.38,
  660.56,
  20.059999999999999,
  0.27,
  3.42,
  0.17,
  7.75,
  4.07,
  42.32,
  97.05,
  70.17,
  2.38,
  5.51,
  0.28999999999999998,
  272.08999999999997,
  11.37,
  66.61,
  17.03,
  104.52,
  213.02,
  27.75,
  3.46,
  217.04,
  19.239999999999998
)

$lastRow = 1 + $values.Length   # header row + one row per value -> A31

for ($i = 0; $i -lt $values.Length; $i++) {
    $cell = $ws.Cells.Item($i + 2, 1)
    $cell.Value = $values[$i]
    # the previous run applied an explicit number-format style (s="1") to
    # every data cell; the fix drops that so cells fall back to the
    # workbook's default (unstyled) formatting.
    $cell.Style = "Normal"
}

$ws.Range("A2:A$lastRow").Select()
